$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the duplicate shared-string key issue: D1, E1, F1 previously all pointed
# at the same "ElementName" string. Give each header its own distinct text.
$ws.Range("D1").Value = "ElementName1"
$ws.Range("E1").Value = "ElementName2"
$ws.Range("F1").Value = "ElementName3"

# Update the view selection/scroll position as recorded in the saved file.
$ws.Range("E6").Select()
$excel.ActiveWindow.ScrollColumn = 3
